# PARSER WORKING - To fix - standard unit
# Restructure the portfolio sheet: split the old "Type" header group into
# several new yield-related columns and move Type/Scheme/AmcName to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for three new columns by shifting the existing H:J
#    (Type, Scheme, AmcName) block to K:M. This preserves their values,
#    formatting and the header row style automatically.
$ws.Range("H1:J1").EntireColumn.Insert()

# 2) Rewrite the header row (A1:M1) with the new lower-cased / renamed
#    column headers. K1/L1/M1 already hold "Type"/"Scheme"/"AmcName" from
#    the shift above, so they are left untouched.
$ws.Range("A1").Value = "name of instrument"
$ws.Range("B1").Value = "isin"
$ws.Range("C1").Value = "coupon"
$ws.Range("D1").Value = "industry"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "market value (mkt)"
$ws.Range("G1").Value = "% to net assets (nav)"
$ws.Range("H1").Value = "yield"
$ws.Range("I1").Value = "yield to call (ytc)"
$ws.Range("J1").Value = "yield to maturity (ytm)"

# 3) The old "Rating/Industry" values in column C are not valid for the new
#    "coupon" column, so they are cleared out.
$lastRow = 22
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = ""
}

# 4) Columns H (yield) and I (yield to call) have no source data, so they
#    stay blank for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = ""
    $ws.Cells.Item($r, 9).Value = ""
}

# 5) Column J becomes "yield to maturity (ytm)" populated with the newly
#    computed per-instrument YTM figures (row 22 - the CDMDF line - has none).
#    Every other text-like column in this sheet stores numbers as plain text,
#    so format the range as Text first to stop Excel auto-converting these
#    values to numbers on assignment.
$ytm = @{
    2  = "0.0763"
    3  = "0.07625"
    4  = "0.077701"
    5  = "0.08554"
    6  = "0.08195"
    7  = "0.0839"
    8  = "0.0813"
    9  = "0.075123"
    10 = "0.081225"
    11 = "0.104622"
    12 = "0.067116"
    13 = "0.0877"
    14 = "0.0835"
    15 = "0.0759"
    16 = "0.075848"
    17 = "0.0759"
    18 = "0.075699"
    19 = "0.07605"
    20 = "0.081499"
    21 = "0.088013"
}
$ws.Range("J2:J" + $lastRow).NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    if ($ytm.ContainsKey($r)) {
        $ws.Cells.Item($r, 10).Value = $ytm[$r]
    } else {
        $ws.Cells.Item($r, 10).Value = ""
    }
}

# 6) Column K ("Type") is rewritten with the normalised instrument-type
#    string (trailing pseudo-nan placeholders come from the upstream parser).
$types = @{
    2  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    3  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    4  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    5  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    6  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    7  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    8  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    9  = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    10 = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    11 = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    12 = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    13 = "Debt Instruments  NAN nan nan nan nan nan nan nan"
    14 = " Securitised Debt  NAN nan nan nan nan nan nan nan"
    15 = "Certificate of Deposit  NAN nan nan nan nan nan nan nan"
    16 = "Certificate of Deposit  NAN nan nan nan nan nan nan nan"
    17 = "Certificate of Deposit  NAN nan nan nan nan nan nan nan"
    18 = "Certificate of Deposit  NAN nan nan nan nan nan nan nan"
    19 = "Commercial Paper  NAN nan nan nan nan nan nan nan"
    20 = "Commercial Paper  NAN nan nan nan nan nan nan nan"
    21 = "Commercial Paper  NAN nan nan nan nan nan nan nan"
    22 = "Corporate Debt Market Development Fund  NAN nan nan nan nan nan nan nan"
}
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 11).Value = $types[$r]
}

Write-Host "Sheet restructured to A1:M22"
